# Append three new data rows (292-294) of FX_IDC:USDUAH quotes to Sheet1,
# extending the existing historical price table by three more months.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: datetime (Excel serial), symbol, open, high, low, close, volume
$newRows = @(
    @(45047.33333333334, "FX_IDC:USDUAH", 36.5684, 36.5685, 36.5681, 36.5681, 0),
    @(45078.33333333334, "FX_IDC:USDUAH", 36.5681, 36.5681, 36.565,  36.565,  0),
    @(45110.33333333334, "FX_IDC:USDUAH", 36.565,  36.565,  36.565,  36.565,  0)
)

$lastExistingRow = 291
$startRow = $lastExistingRow + 1

# Copy the formatting (number format/alignment/borders/font) of the last
# existing date cell so the newly appended date cells keep the same look.
$ws.Cells.Item($lastExistingRow, 1).Copy() | Out-Null

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
}

$excel.CutCopyMode = 0
